$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "S_4_withindex_sequence_TGAGGTT.fastq.gz"
$ws.Range("F3").Value = "S_4_withindex_sequence_GCTTAGA.fastq.gz"
$ws.Range("F4").Value = "s_4_withindex_sequence_ATGACAG.fastq.gz"
$ws.Range("F5").Value = "s_4_withindex_sequence_CACCTCC.fastq.gz"
$ws.Range("F6").Value = "s_4_withindex_sequence_ATCGAGC.fastq.gz"
$ws.Range("F7").Value = "s_4_withindex_sequence_TACTCTA.fastq.gz"
$ws.Range("F8").Value = "S_4_withindex_sequence_AGACTGA.fastq.gz"
$ws.Range("F9").Value = "S_4_withindex_sequence_CTTGGAA.fastq.gz"
$ws.Range("F10").Value = "s_4_withindex_sequence_CCGATTA.fastq.gz"
$ws.Range("F11").Value = "s_4_withindex_sequence_GGCAGCG.fastq.gz"
$ws.Range("F12").Value = "s_4_withindex_sequence_CCATCAT.fastq.gz"
$ws.Range("F13").Value = "s_4_withindex_sequence_TAACAAG.fastq.gz"
$ws.Range("F14").Value = "s_4_withindex_sequence_GAGGCGT.fastq.gz"
$ws.Range("F15").Value = "s_4_withindex_sequence_TTTAACT.fastq.gz"
$ws.Range("F16").Value = "s_4_withindex_sequence_GGTCCTC.fastq.gz"
$ws.Range("F17").Value = "s_4_withindex_sequence_CGGTGGC.fastq.gz"
$ws.Range("F18").Value = "s_4_withindex_sequence_ACTGTCG.fastq.gz"
$ws.Range("F19").Value = "s_4_withindex_sequence_GTATTTG.fastq.gz"

$ws.Range("F20").Select()

$ws.Rows.Item(1).RowHeight = 13.8

